$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")

# Append the new row (66) with the August 4th SSA raw/clean data, matching
# the existing table layout (Fecha, Confirmados, Negativos, Sospechosos,
# Defunciones, Porcentaje hospitalizados).
#
# Column A stores the date as plain text (e.g. "2020-08-03" in row 65), so
# force a text number format before assigning the string value - otherwise
# Excel's smart entry would silently reinterpret "2020-08-04" as a date
# serial. Clear the format again afterwards so the cell ends up with the
# same (default) style as its neighbours in column A.
$dateCell = $ws.Cells.Item(66, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2020-08-04"
$dateCell.ClearFormats()

$ws.Cells.Item(66, 2).Value = 449961
$ws.Cells.Item(66, 3).Value = 493873
$ws.Cells.Item(66, 4).Value = 82460
$ws.Cells.Item(66, 5).Value = 48869
$ws.Cells.Item(66, 6).Value = 26.86
